# "User yandex api in registration form"
# Adds two new registration entries (Name / Phone) to the bottom of the
# active sheet ("Приморский" — sheet2), which grows the used range from
# A1:B7 to A1:B9, appends 4 new shared strings, widens column B slightly
# to fit the new phone numbers, and moves the selection to the new last
# cell (B9) — matching the recorded OOXML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 8 and 9: Имя / Телефон
$ws.Range("A8").Value = "Алексей"
$ws.Range("B8").Value = "8-905-223-02-49"
$ws.Range("A9").Value = "Борис"
$ws.Range("B9").Value = "8-962-685-01-80"

# Column B grows from 14.28515625 to 15 (OOXML "width" units) to fit the
# new, slightly longer phone numbers.
$ws.Columns.Item(2).ColumnWidth = 14.2

# Active cell / selection follows the newly entered last row.
$ws.Range("B9").Select() | Out-Null
